# refactor currency conversion, now explicit source and target amounts
#
# currency_conversions sheet ("foreign_amount" column) is split into
# explicit source_amount / source_fees / source_currency and
# target_amount / target_fees / target_currency columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Old layout:  A date | B foreign_amount | C source_fees | D source_currency | E target_currency | F comment
# New layout:  A date | B source_amount  | C source_fees | D source_currency | E target_amount | F target_fees | G target_currency | H comment
#
# Insert two fresh columns right before the old "target_currency" column (E)
# so the row grows from 6 to 8 columns: one becomes target_amount (E),
# the other becomes target_fees (F); old target_currency/comment shift to G/H.
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).Insert()

# Header row (target_amount is assigned first so new shared strings keep the
# same ordering as the authoritative edit: target_amount, source_amount,
# target_fees)
$ws.Range("E1").Value2 = "target_amount"

$ws.Range("B1").Value2 = "source_amount"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.ThemeColor = 1

$ws.Range("F1").Value2 = "target_fees"

# Data row: the old "foreign_amount" (144.74) becomes the target_amount (money
# received in the target currency); source_amount is now explicit (-1 unit of
# source currency spent), and target_fees defaults to 0.
$ws.Range("B2").Value2 = -1
$ws.Range("E2").Value2 = 144.74
$ws.Range("F2").Value2 = 0

# This sheet becomes the active / selected one.
$ws.Activate()
